# Bottle real pouring results
# Adds 4 new scene rows (real-world pouring measurements) to the summary sheet,
# widens the "path" column so the longer paths are readable, and leaves the
# trailing rows of the newly selected block formatted-but-empty (Text format),
# matching how the data was pasted/typed into the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The user selected A52:L63 and set it to Text format before typing in the
# new rows, so every cell in that block (including the still-empty rows
# 56:63) keeps the "Text" number format.
$ws.Range("A52:L63").NumberFormat = "@"

# New row 52
$ws.Range("B52").Value = "../../output/MediumBottle/Medium_370_1399_78"
$ws.Range("C52").Value = "0.03"
$ws.Range("D52").Value = "78.0"
$ws.Range("E52").Value = "1.4"
$ws.Range("F52").Value = "369.9"
$ws.Range("G52").Value = "273.10"
$ws.Range("H52").Value = "249.60"
$ws.Range("I52").Value = "23.50"

# New row 53
$ws.Range("B53").Value = "../../output/MediumBottle/Medium_350_1399_54"
$ws.Range("C53").Value = "0.03"
$ws.Range("D53").Value = "54.0"
$ws.Range("E53").Value = "1.4"
$ws.Range("F53").Value = "349.9"
$ws.Range("G53").Value = "87.3"
$ws.Range("H53").Value = "78.5"
$ws.Range("I53").Value = "8.80"

# New row 54
$ws.Range("B54").Value = "../../output/MediumBottle/Medium_430_1000_46"
$ws.Range("C54").Value = "0.03"
$ws.Range("D54").Value = "46.0"
$ws.Range("E54").Value = "1.0"
$ws.Range("F54").Value = "429.89"
$ws.Range("G54").Value = "87.03"
$ws.Range("H54").Value = "79.04"
$ws.Range("I54").Value = "7.99"

# New row 55
$ws.Range("B55").Value = "../../output/MediumBottle/Medium_410_1399_90"
$ws.Range("C55").Value = "0.03"
$ws.Range("D55").Value = "90.0"
$ws.Range("E55").Value = "1.4"
$ws.Range("F55").Value = "409.892"
$ws.Range("G55").Value = "380.678"
$ws.Range("H55").Value = "341.74"
$ws.Range("I55").Value = "38.935"

# volume_start (col F) carries the same yellow highlight used throughout the
# rest of the sheet for that column.
$ws.Range("F52:F55").Interior.Color = 65535

# Widen the path column (B) now that it holds longer relative paths; the
# rest of the columns keep their original width.
$ws.Columns.Item(2).ColumnWidth = 44.33

# Scroll/select near the newly entered data, like the author did when done.
$ws.Range("I60").Select()
